$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fix the 2018.10.15 week's "er图" task row (row 47): shorten the task
# description and bump the completion percentage from 50% to 90%. ---
$ws.Range("B47").Value = "初步设计完成er图"
$ws.Range("C47").Value = "未完成（90%）"

# --- Append a brand-new weekly block (2018.10.17) mirroring the existing
# layout used for the other week tables. Copy the formatting from the most
# recent block (rows 42:50, i.e. header + 5 members + summary footer) down
# to rows 52:60, leaving row 51 blank as a spacer (matches the existing
# pattern at row 41). ---
$ws.Range("A42:D50").Copy()
$ws.Range("A52").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 52: week date header (merged A52:D52)
$ws.Range("A52").Value = "日期：2018.10.17 第七周周三"
$ws.Range("A52:D52").Merge()

# Row 53: column headers
$ws.Range("A53").Value = "组员"
$ws.Range("B53").Value = "计划内容"
$ws.Range("C53").Value = "完成情况"
$ws.Range("D53").Value = "备注"

# Rows 54-58: member tasks
$ws.Range("A54").Value = "邱志鹏"
$ws.Range("B54").Value = "再一次审查用例分析文档中的内容，将不足之处完善"
$ws.Range("C54").Value = "未完成"

$ws.Range("A55").Value = "黄立根"
$ws.Range("B55").Value = "开始编写网页端管理员界面网页设计（仅html+css）"
$ws.Range("C55").Value = "未完成"

$ws.Range("A56").Value = "黄俊贤"
$ws.Range("B56").Value = "设计完成数据库管理员表、群成员表、群管理表、聊天信息表设计"
$ws.Range("C56").Value = "未完成"

$ws.Range("A57").Value = "李达波"
$ws.Range("B57").Value = "完成数据库er图设计，编写网络交互API文档"
$ws.Range("C57").Value = "未未完成"

$ws.Range("A58").Value = "冯德志"
$ws.Range("B58").Value = "设计完成数据库用户表、群组表、好友表设计"
$ws.Range("C58").Value = "未完成"

# Rows 59-60: summary footer (merged A59:D60)
$ws.Range("A59").Value = "总结："
$ws.Range("A59:D60").Merge()

# Match the author's final selection/view position after the edit.
$ws.Range("B54").Select()
